$wb = $excel.ActiveWorkbook

# Both sheets ("Time mismatch" and "Midnight samples") share the same
# column layout: B = serial_num, ... M = tds_ppt (to be removed), with
# ammonia_mg_per_L/follow_up/month_year/hour/time_mismatch shifting left.
foreach ($ws in $wb.Worksheets) {
    # Rename the serial_num header to sample_idx.
    $ws.Range("B1").Value = "sample_idx"

    # Delete the whole tds_ppt column (M); everything to its right
    # (ammonia_mg_per_L, follow_up, month_year, hour, time_mismatch)
    # shifts one column to the left automatically.
    $ws.Range("M1").EntireColumn.Delete()
}
